$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 8666.333000000001
$ws.Cells.Item(62, 9).Value = 8499.5
$ws.Cells.Item(62, 11).Value = 8499.5
$ws.Cells.Item(62, 13).Value = -7875.5
$ws.Cells.Item(65, 8).Value = 8666.333000000001
$ws.Cells.Item(65, 9).Value = 8499.5
$ws.Cells.Item(65, 11).Value = 42497.5
$ws.Cells.Item(65, 13).Value = -39377.5
$ws.Cells.Item(80, 8).Value = 464.57144
$ws.Cells.Item(80, 9).Value = 410
$ws.Cells.Item(80, 10).Value = 505.5
$ws.Cells.Item(80, 11).Value = 1230
$ws.Cells.Item(80, 12).Value = 1516.5
$ws.Cells.Item(80, 13).Value = -232
$ws.Cells.Item(80, 14).Value = -3512.5
$ws.Cells.Item(83, 8).Value = 464.57144
$ws.Cells.Item(83, 9).Value = 410
$ws.Cells.Item(83, 10).Value = 505.5
$ws.Cells.Item(83, 11).Value = 3690
$ws.Cells.Item(83, 12).Value = 4549.5
$ws.Cells.Item(83, 13).Value = 1302
$ws.Cells.Item(83, 14).Value = -14533.5
$ws.Cells.Item(113, 8).Value = 19629
$ws.Cells.Item(113, 10).Value = 19629
$ws.Cells.Item(113, 12).Value = 19629
$ws.Cells.Item(113, 14).Value = -26137
$ws.Cells.Item(135, 8).Value = 1166.3334
$ws.Cells.Item(135, 9).Value = 1368.2858
$ws.Cells.Item(135, 10).Value = 459.5
$ws.Cells.Item(135, 11).Value = 12314.5722
$ws.Cells.Item(135, 12).Value = 4135.5
$ws.Cells.Item(135, 13).Value = -9779.572200000001
$ws.Cells.Item(135, 14).Value = -9205.5
$ws.Cells.Item(141, 8).Value = 41415.832
$ws.Cells.Item(141, 9).Value = 41415.832
$ws.Cells.Item(141, 11).Value = 124247.496
$ws.Cells.Item(141, 13).Value = -119067.496

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5054.6665
$ws.Cells.Item(32, 9).Value = 3746.64
$ws.Cells.Item(32, 10).Value = 11594.8
$ws.Cells.Item(32, 11).Value = 3746.64
$ws.Cells.Item(32, 12).Value = 11594.8
$ws.Cells.Item(32, 13).Value = -3459.64
$ws.Cells.Item(32, 14).Value = -12168.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 665
$ws.Cells.Item(80, 9).Value = 528.1429000000001
$ws.Cells.Item(80, 11).Value = 528.1429000000001
$ws.Cells.Item(80, 13).Value = 469.8570999999999
$ws.Cells.Item(83, 8).Value = 665
$ws.Cells.Item(83, 9).Value = 528.1429000000001
$ws.Cells.Item(83, 11).Value = 2640.7145
$ws.Cells.Item(83, 13).Value = 2351.2855
$ws.Cells.Item(105, 8).Value = 4667.5
$ws.Cells.Item(105, 9).Value = 4667.5
$ws.Cells.Item(105, 11).Value = 4667.5
$ws.Cells.Item(105, 13).Value = -2920.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1383.1666
$ws.Cells.Item(16, 9).Value = 1420
$ws.Cells.Item(16, 10).Value = 1199
$ws.Cells.Item(16, 11).Value = 1420
$ws.Cells.Item(16, 12).Value = 1199
$ws.Cells.Item(16, 13).Value = -1133
$ws.Cells.Item(16, 14).Value = -1773
$ws.Cells.Item(22, 8).Value = 689
$ws.Cells.Item(22, 9).Value = 823.75
$ws.Cells.Item(22, 11).Value = 823.75
$ws.Cells.Item(22, 13).Value = -473.75
$ws.Cells.Item(32, 8).Value = 4999.5
$ws.Cells.Item(32, 9).Value = 4999.5
$ws.Cells.Item(32, 11).Value = 4999.5
$ws.Cells.Item(32, 13).Value = -4683.5
$ws.Cells.Item(58, 8).Value = 2370.375
$ws.Cells.Item(58, 9).Value = 2327.1667
$ws.Cells.Item(58, 11).Value = 2327.1667
$ws.Cells.Item(58, 13).Value = -2124.1667
$ws.Cells.Item(99, 8).Value = 1716
$ws.Cells.Item(99, 9).Value = 1574
$ws.Cells.Item(99, 10).Value = 2000
$ws.Cells.Item(99, 11).Value = 1574
$ws.Cells.Item(99, 12).Value = 2000
$ws.Cells.Item(99, 13).Value = -76
$ws.Cells.Item(99, 14).Value = -4996
$ws.Cells.Item(113, 8).Value = 1383.1666
$ws.Cells.Item(113, 9).Value = 1420
$ws.Cells.Item(113, 10).Value = 1199
$ws.Cells.Item(113, 11).Value = 1420
$ws.Cells.Item(113, 12).Value = 1199
$ws.Cells.Item(113, 13).Value = 750
$ws.Cells.Item(113, 14).Value = -5539
$ws.Cells.Item(126, 8).Value = 1716
$ws.Cells.Item(126, 9).Value = 1574
$ws.Cells.Item(126, 10).Value = 2000
$ws.Cells.Item(126, 11).Value = 4722
$ws.Cells.Item(126, 12).Value = 6000
$ws.Cells.Item(126, 13).Value = -2252
$ws.Cells.Item(126, 14).Value = -10940
$ws.Cells.Item(134, 8).Value = 2375
$ws.Cells.Item(134, 9).Value = 2375
$ws.Cells.Item(134, 11).Value = 7125
$ws.Cells.Item(134, 13).Value = -4590
$ws.Cells.Item(136, 8).Value = 2370.375
$ws.Cells.Item(136, 9).Value = 2327.1667
$ws.Cells.Item(136, 11).Value = 6981.500100000001
$ws.Cells.Item(136, 13).Value = -4431.500100000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 11.416667
$ws.Cells.Item(2, 9).Value = 13.25
$ws.Cells.Item(2, 10).Value = 7.75
$ws.Cells.Item(2, 11).Value = 79.5
$ws.Cells.Item(2, 12).Value = 46.5
$ws.Cells.Item(2, 13).Value = 33.5
$ws.Cells.Item(2, 14).Value = -272.5
$ws.Cells.Item(68, 8).Value = 1966.3334
$ws.Cells.Item(68, 10).Value = 1949.5
$ws.Cells.Item(68, 12).Value = 5848.5
$ws.Cells.Item(68, 14).Value = -7470.5
$ws.Cells.Item(71, 8).Value = 1966.3334
$ws.Cells.Item(71, 10).Value = 1949.5
$ws.Cells.Item(71, 12).Value = 17545.5
$ws.Cells.Item(71, 14).Value = -25657.5
$ws.Cells.Item(113, 8).Value = 719.25
$ws.Cells.Item(113, 10).Value = 812
$ws.Cells.Item(113, 12).Value = 2436
$ws.Cells.Item(113, 14).Value = -6776

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29, 8).Value = 4500
$ws.Cells.Item(29, 10).Value = 4500
$ws.Cells.Item(29, 12).Value = 4500
$ws.Cells.Item(29, 14).Value = -5080
$ws.Cells.Item(35, 8).Value = 15
$ws.Cells.Item(35, 9).Value = 15
$ws.Cells.Item(35, 11).Value = 15
$ws.Cells.Item(35, 13).Value = 283
$ws.Cells.Item(70, 8).Value = 5681.5
$ws.Cells.Item(70, 9).Value = 5681.5
$ws.Cells.Item(70, 11).Value = 5681.5
$ws.Cells.Item(70, 13).Value = -5411.5
$ws.Cells.Item(73, 8).Value = 5681.5
$ws.Cells.Item(73, 9).Value = 5681.5
$ws.Cells.Item(73, 11).Value = 5681.5
$ws.Cells.Item(73, 13).Value = -4745.5
$ws.Cells.Item(80, 8).Value = 2346
$ws.Cells.Item(80, 10).Value = 2640
$ws.Cells.Item(80, 12).Value = 2640
$ws.Cells.Item(80, 14).Value = -4636
$ws.Cells.Item(83, 8).Value = 2346
$ws.Cells.Item(83, 10).Value = 2640
$ws.Cells.Item(83, 12).Value = 13200
$ws.Cells.Item(83, 14).Value = -23184

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(33, 8).Value = 15
$ws.Cells.Item(33, 9).Value = 15
$ws.Cells.Item(33, 11).Value = 15
$ws.Cells.Item(33, 13).Value = 275
$ws.Cells.Item(40, 8).Value = 3097.6
$ws.Cells.Item(40, 9).Value = 3442.25
$ws.Cells.Item(40, 10).Value = 1719
$ws.Cells.Item(40, 11).Value = 3442.25
$ws.Cells.Item(40, 12).Value = 1719
$ws.Cells.Item(40, 13).Value = -3306.25
$ws.Cells.Item(40, 14).Value = -1991

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3520.2
$ws.Cells.Item(136, 9).Value = 3520.2
$ws.Cells.Item(136, 11).Value = 10560.6
$ws.Cells.Item(136, 13).Value = -8010.599999999999

Write-Host "Applied Kraken_Profits cell updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
